# Generate Report for Handoff
#
# The localization status report moves from "In Translation" to
# "Ready for handoff": the Status cells are updated and the handoff/
# generate timestamps are bumped to the moment the handoff report was
# produced. Column widths for the (now wider) status/date text are
# widened to fit the new content on all three sheets.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Status columns for zh-cn (E) and de-de (F)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2017-02-17 09:20:59"

# Widen the Status columns to fit "Ready for handoff"
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# ---- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2017-02-17 09:20:42"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# ---- de-de sheet ---------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2017-02-17 09:20:59"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.33

Write-Host "Report regenerated for handoff"
